$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.998.98"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "2.981.17"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Formula = '="596.45"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Formula = '="145.98"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "2.979.70"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D11").Formula = '="0.144"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +2.51%  "
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("D13").Formula = '="0.0000236"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +5.15%  "
$ws.Range("D14").Formula = '="33.68"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "3.474.34"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "62.837.24"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "2.997.25"
$ws.Range("E19").Value = "  +2.85%  "
$ws.Range("D20").Formula = '="444.58"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").Formula = '="13.61"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("D22").Formula = '="0.678"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").Formula = '="82.34"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").Formula = '="11.12"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("D26").Formula = '="12.18"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").Formula = '="7.17"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("E31").Value = "  -5.84%  "
$ws.Range("D32").Formula = '="26.66"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "0.0₃0886"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").Formula = '="0.995"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").Formula = '="5.68"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").Formula = '="2.06"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +4.05%  "
$ws.Range("D39").Formula = '="49.83"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Formula = '="2.98"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").Formula = '="8.65"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Formula = '="38.98"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -6.66%  "
$ws.Range("D45").Formula = '="375.58"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.712.52"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Formula = '="0.0345"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Formula = '="135.01"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Formula = '="23.46"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Formula = '="0.105"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.01%  "

$excel.CutCopyMode = 0

